$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first data bucket row "(1925, 1935]" / 229 -> all subsequent
# rows shift up by one, and the sheet's used range shrinks from B11 to B10.
$ws.Rows.Item(2).Delete()
